# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell "D2" "41.537.48"
Set-TextCell "E2" "  +0.06%  "

Set-TextCell "D3" "2.454.64"
Set-TextCell "E3" "  -1.22%  "

Set-TextCell "E4" "  -0.12%  "

Set-TextCell "D5" "314.40"

Set-TextCell "D6" "92.02"
Set-TextCell "E6" "  +0.14%  "

Set-TextCell "D7" "0.546"
Set-TextCell "E7" "  +1.76%  "

Set-TextCell "E8" "  -0.20%  "

Set-TextCell "D9" "0.507"
Set-TextCell "E9" "  +3.48%  "

Set-TextCell "D10" "32.22"
Set-TextCell "E10" "  +0.48%  "

Set-TextCell "E11" "  +2.91%  "

Set-TextCell "E12" "  +0.85%  "

Set-TextCell "D13" "2.836.73"
Set-TextCell "E13" "  -1.12%  "

Set-TextCell "D14" "6.80"
Set-TextCell "E14" "  +0.47%  "

Set-TextCell "D15" "15.75"
Set-TextCell "E15" "  +2.89%  "

Set-TextCell "D16" "2.466.56"
Set-TextCell "E16" "  +0.15%  "

Set-TextCell "D17" "0.772"
Set-TextCell "E17" "  +1.77%  "

Set-TextCell "D18" "41.564.46"
Set-TextCell "E18" "  +0.23%  "

Set-TextCell "D19" "6.45"
Set-TextCell "E19" "  +2.51%  "

Set-TextCell "D20" "0.0₃0933"
Set-TextCell "E20" "  +1.81%  "

Set-TextCell "D21" "70.67"
Set-TextCell "E21" "  +0.50%  "

Set-TextCell "D22" "11.30"
Set-TextCell "E22" "  +2.21%  "

Set-TextCell "D23" "237.54"
Set-TextCell "E23" "  +1.53%  "

Set-TextCell "E24" "  +0.48%  "

Set-TextCell "E25" "  -0.10%  "

Set-TextCell "D26" "1.89"
Set-TextCell "E26" "  +0.79%  "

Set-TextCell "D27" "24.22"
Set-TextCell "E27" "  +0.03%  "

Set-TextCell "E28" "  +0.71%  "

Set-TextCell "D29" "9.64"
Set-TextCell "E29" "  +0.77%  "

Set-TextCell "D30" "34.86"
Set-TextCell "E30" "  -3.75%  "

Set-TextCell "D31" "155.73"
Set-TextCell "E31" "  +1.51%  "

Set-TextCell "D32" "5.44"
Set-TextCell "E32" "  +1.65%  "

Set-TextCell "E33" "  +0.46%  "

Set-TextCell "D34" "0.0756"

Set-TextCell "D35" "2.47"
Set-TextCell "E35" "  -0.68%  "

Set-TextCell "D36" "17.39"
Set-TextCell "E36" "  -4.50%  "

Set-TextCell "D37" "2.87"
Set-TextCell "E37" "  -3.40%  "

Set-TextCell "E38" "  +1.73%  "

Set-TextCell "E39" "  +1.68%  "

Set-TextCell "D40" "1.78"
Set-TextCell "E40" "  -2.23%  "

Set-TextCell "D41" "3.93"
Set-TextCell "E41" "  -2.41%  "

Set-TextCell "E42" "  -0.39%  "

Set-TextCell "D43" "1.965.41"
Set-TextCell "E43" "  +0.99%  "

Set-TextCell "B44" "EnergySwap"
Set-TextCell "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D44" "18.73"
Set-TextCell "E44" "  -4.29%  "

Set-TextCell "B45" "VeChain"
Set-TextCell "C45" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D45" "0.0280"
Set-TextCell "E45" "  +0.19%  "

Set-TextCell "D46" "2.89"
Set-TextCell "E46" "  -1.44%  "

Set-TextCell "E47" "  +2.44%  "

Set-TextCell "D48" "2.694.88"
Set-TextCell "E48" "  -0.71%  "

Set-TextCell "D49" "96.16"
Set-TextCell "E49" "  +0.74%  "

Set-TextCell "D50" "66.25"
Set-TextCell "E50" "  +0.06%  "

Set-TextCell "E51" "  -2.05%  "
